$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse the distinct per-row IP addresses (192.168.1.113 .. 192.168.1.117)
# down to a single reused "127.0.0.1" value for rows 2-6.
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("F3").Value = "127.0.0.1"
$ws.Range("F4").Value = "127.0.0.1"
$ws.Range("F5").Value = "127.0.0.1"
$ws.Range("F6").Value = "127.0.0.1"

# Move the active selection, as it was left after editing.
$ws.Range("F14").Select()
